# Apply strikethrough formatting to the requirement bullet points that are
# being merged/superseded by the new "input feature" (per commit message
# "changes to merge the input feature"). The five bullets struck through are
# the original map/marker related requirements.

$d = $word.ActiveDocument

$targets = @(
    "Focus is on Map.",
    "At least 10 markers should be on the Map.",
    "When the marker is clicked, the info window should pop up with info about that location. Ex.: Address, name of the location, etc.",
    "There should be 2 buttons to filter markers. When a button is clicked on a certain type of location marker should be on map. Others should disappear. Ex.: a button to show only museums and a button to show only waterfalls",
    "Allow users to put a marker for their location on the Map via entering an address into an input text."
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    foreach ($t in $targets) {
        if ($text.StartsWith($t)) {
            $p.Range.Font.StrikeThrough = 1
            break
        }
    }
}
